$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
"B2=1.02",
"C2=1.059692457733099",
"D2=1.0578066871998",
"E2=1.063386412290673",
"F2=1.066211644921897",
"I2=1.053133242117264",
"J2=1.064677168939965",
"K2=1.060540648049058",
"L2=1.06610519161132",
"M2=1.068922802805314",
"N2=1.066189132640038",
"B3=1.02",
"C3=1.061258192075281",
"D3=1.059047221175747",
"E3=1.064903673854067",
"F3=1.067940693955115",
"I3=1.053768880442332",
"J3=1.065893312102734",
"K3=1.061594008123482",
"L3=1.06743569217258",
"M3=1.07046512228574",
"N3=1.067407002865591",
"B4=1.02",
"C4=1.062267621809329",
"D4=1.05984610034646",
"E4=1.065882088920246",
"F4=1.069056383198932",
"I4=1.054176705916015",
"J4=1.066676150915918",
"K4=1.062271219681235",
"L4=1.068292761771585",
"M4=1.071459528984112",
"N4=1.068190953399683",
"B5=1.02",
"C5=1.062691113654707",
"D5=1.060181044143742",
"E5=1.066292625047255",
"F5=1.069524687310988",
"I5=1.054347330993647",
"J5=1.06700428994103",
"K5=1.062554881278892",
"L5=1.068652164012606",
"M5=1.071876735618198",
"N5=1.06851955841985",
"B6=1.02",
"C6=1.062762169061321",
"D6=1.060237229949565",
"E6=1.066361509959568",
"F6=1.069603275176953",
"I6=1.054375931538193",
"J6=1.067059329604294",
"K6=1.06260244875313",
"L6=1.068712456243309",
"M6=1.071946737492594",
"N6=1.068574676245748",
"B7=1.02",
"C7=1.06227328393614",
"D7=1.059850579421864",
"E7=1.065887577607996",
"F7=1.069062643555016",
"I7=1.05417898904753",
"J7=1.066680539306472",
"K7=1.062275014045354",
"L7=1.068297567679567",
"M7=1.071465107008106",
"N7=1.068195348022255",
"B8=1.02",
"C8=1.06022238097189",
"D8=1.058226731435151",
"E8=1.063899880725087",
"F8=1.066796640787877",
"I8=1.053348782935503",
"J8=1.065089023778132",
"K8=1.060897550244118",
"L8=1.06655564599754",
"M8=1.0694447862943",
"N8=1.066601572359316",
"B9=1.02",
"C9=1.05657939267026",
"D9=1.055335446361488",
"E9=1.060371001158087",
"F9=1.062779048598814",
"I9=1.051858914742795",
"J9=1.062252722882769",
"K9=1.058436213214279",
"L9=1.063456074751248",
"M9=1.065856674922897",
"N9=1.063761243591352",
"B10=1.02",
"C10=1.054130285210104",
"D10=1.053387102214066",
"E10=1.057999869157936",
"F10=1.060083123019009",
"I10=1.050847099912191",
"J10=1.060339676497451",
"K10=1.056771707652847",
"L10=1.061368663358006",
"M10=1.06344482600753",
"N10=1.061845480460739",
"B11=1.02",
"C11=1.053064736586886",
"D11=1.052538341662403",
"E11=1.05696855094001",
"F11=1.058911378095976",
"I11=1.050404465764967",
"J11=1.059505873531476",
"K11=1.056045198498795",
"L11=1.060459622382502",
"M11=1.062395572933076",
"N11=1.061010493398913",
"B12=1.02",
"C12=1.052668164782617",
"D12=1.052222291199217",
"E12=1.05658476519957",
"F12=1.058475460912176",
"I12=1.050239364954146",
"J12=1.059195329525716",
"K12=1.055774460358504",
"L12=1.06012117044218",
"M12=1.062005078620513",
"N12=1.060699508385045",
"B13=1.02",
"C13=1.052753266351754",
"D13=1.052290120765103",
"E13=1.05666712086442",
"F13=1.058568997662435",
"I13=1.050274810861815",
"J13=1.05926198018209",
"K13=1.055832574667641",
"L13=1.060193805700723",
"M13=1.062088875430505",
"N13=1.060766253692996",
"B14=1.02",
"C14=1.053031971846735",
"D14=1.05251223289058",
"E14=1.056936841636567",
"F14=1.058875359007813",
"I14=1.050390832542783",
"J14=1.059480220936303",
"K14=1.056022837237993",
"L14=1.060431662121246",
"M14=1.062363310057254",
"N14=1.060984804374113",
"B15=1.02",
"C15=1.05320358765179",
"D15=1.05264897930673",
"E15=1.057102931152692",
"F15=1.059064027798357",
"I15=1.050462226054622",
"J15=1.059614575507692",
"K15=1.056139947223006",
"L15=1.060578107673068",
"M15=1.062532297705354",
"N15=1.061119349744404",
"B16=1.02",
"C16=1.054200893830745",
"D16=1.05344332265596",
"E16=1.058068215957064",
"F16=1.060160793716725",
"I16=1.050876380265751",
"J16=1.060394897339709",
"K16=1.056819800987022",
"L16=1.061428882982237",
"M16=1.063514356584504",
"N16=1.061900779722926",
"B17=1.02",
"C17=1.054825107879148",
"D17=1.053940213674116",
"E17=1.058672469640804",
"F17=1.060847576834554",
"I17=1.051134954078406",
"J17=1.060882905321287",
"K17=1.057244701526882",
"L17=1.061961154520042",
"M17=1.064129050016581",
"N17=1.062389480731813",
"B18=1.02",
"C18=1.055188713525416",
"D18=1.054229549300619",
"E18=1.059024477669151",
"F18=1.061247743414337",
"I18=1.051285341110977",
"J18=1.061167028139328",
"K18=1.057491983309022",
"L18=1.062271120965229",
"M18=1.064487118098905",
"N18=1.062674007036845",
"B19=1.02",
"C19=1.055312611449756",
"D19=1.054328122174262",
"E19=1.059144428687073",
"F19=1.061384118808734",
"I19=1.051336545739649",
"J19=1.061263818251643",
"K19=1.057576206197217",
"L19=1.062376727424155",
"M19=1.06460913050541",
"N19=1.062770934602222",
"B20=1.02",
"C20=1.054758186209781",
"D20=1.053886953016891",
"E20=1.058607684873639",
"F20=1.060773935334229",
"I20=1.051107256582822",
"J20=1.060830600989122",
"K20=1.057199171279299",
"L20=1.061904098500703",
"M20=1.064063148198942",
"N20=1.062337102121497",
"B21=1.02",
"C21=1.052949921691971",
"D21=1.052446848138392",
"E21=1.056857435278316",
"F21=1.058785162125093",
"I21=1.05035668608051",
"J21=1.059415977574844",
"K21=1.055966834065514",
"L21=1.06036164136946",
"M21=1.062282516885988",
"N21=1.060920469779711",
"B22=1.02",
"C22=1.05180847437349",
"D22=1.051536859467953",
"E22=1.055752877796587",
"F22=1.057530803609412",
"I22=1.04988079459499",
"J22=1.05852172377937",
"K22=1.055186914666955",
"L22=1.059387237629913",
"M22=1.061158586089626",
"N22=1.060024946041276",
"B23=1.02",
"C23=1.052414011690841",
"D23=1.052019696706629",
"E23=1.056338819482908",
"F23=1.058196142527448",
"I23=1.050133453766644",
"J23=1.058996246935748",
"K23=1.05560085278812",
"L23=1.059904228783472",
"M23=1.061754823996316",
"N23=1.060500143074965",
"B24=1.02",
"C24=1.054788426738663",
"D24=1.053911020736869",
"E24=1.058636959689635",
"F24=1.060807212063201",
"I24=1.051119773231875",
"J24=1.06085423668422",
"K24=1.057219746154603",
"L24=1.061929881200081",
"M24=1.064092927852364",
"N24=1.062360771381994",
"B25=1.02",
"C25=1.057524728024344",
"D25=1.056086527367964",
"E25=1.061286506170715",
"F25=1.063820709301783",
"I25=1.052247320486008",
"J25=1.06298982706278",
"K25=1.059076637579592",
"L25=1.064261035389568",
"M25=1.066787705255017",
"N25=1.064499394543847"
)

foreach ($entry in $updates) {
    $parts = $entry.Split("=")
    $cellRef = $parts[0]
    $newValue = [double]$parts[1]
    $ws.Range($cellRef).Value = $newValue
}

Write-Output "Updated $($updates.Count) cells"
